$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.142.90'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.599.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '302.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3781'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '52.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3607'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.261'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.001'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08110'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.559'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.374'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001244'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.601.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06880'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.525'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.150.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.381'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.977'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.253'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.387'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.769'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.778.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9728'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07481'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02700'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2503'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.54%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.08798'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.057'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.360'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7069'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6510'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.304'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07951'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.201'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.216'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.39%  '
